$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -----------------------------------------------------------------
# 1. Insert a blank row at row 3 (old rows 3,4,5 shift down to 4,5,6)
#    Clear it completely so it stays a truly empty row.
# -----------------------------------------------------------------
$ws.Rows.Item(3).Insert()
$ws.Rows.Item(3).Clear()

# Drop the stray empty A2:C2 cells left over from the old layout.
$ws.Range("A2:C2").Clear()

# -----------------------------------------------------------------
# 2. Header row 1: D1:E1 merged "Lines of Code", F1:G1 merged "Tomorrow"
# -----------------------------------------------------------------
$ws.Range("D1").Value = "Lines of Code"
$ws.Range("F1").Value = "Tomorrow"
$ws.Range("D1:E1").Merge()
$ws.Range("F1:G1").Merge()

$ws.Range("A1:G1").Font.Bold = $true
$ws.Range("A1:G1").HorizontalAlignment = -4108
$ws.Range("A1:G1").VerticalAlignment = -4108
$ws.Range("A1:G1").WrapText = $true

# -----------------------------------------------------------------
# 3. Sub-header row 2: D2 "+", E2 "-", F2 "Primary Goals", G2 "Secondary Goals"
# -----------------------------------------------------------------
$ws.Range("D2").Value = "+"
$ws.Range("E2").Value = "-"
$ws.Range("F2").Value = "Primary Goals"
$ws.Range("G2").Value = "Secondary Goals"

$ws.Range("D2:G2").Font.Bold = $true
$ws.Range("D2:F2").HorizontalAlignment = -4108
$ws.Range("D2:G2").VerticalAlignment = -4108
$ws.Range("D2:G2").WrapText = $true

# -----------------------------------------------------------------
# 4. Data row 4 (was row 3): day 1
# -----------------------------------------------------------------
$ws.Range("C4").Value = "Studied the process of crawling data from any webpage, watched python tutorials, built a basic web-scraper by following a youtube tutorial. Reviewed Early Warning Systems (EWS) and understood the flow of the model."
$ws.Range("D4").Value = 210
$ws.Range("E4").Value = 0
$ws.Range("F4").Value = "To Build a basic web scraper that scrapes data of one company from one website and start working from there."
$ws.Range("G4").ClearContents()

$ws.Range("A4:F4").HorizontalAlignment = -4108
$ws.Range("A4:G4").VerticalAlignment = -4108
$ws.Range("A4:G4").WrapText = $true
$ws.Range("F4").Interior.Color = 5296274

# -----------------------------------------------------------------
# 5. Data row 5 (was row 4): day 2
# -----------------------------------------------------------------
$ws.Range("D5").Value = 885
$ws.Range("E5").Value = 170
$ws.Range("F5").Value = "Clean article data in csv file. Include more companies and a new website."
$ws.Range("G5").Value = "Build keywords dictionary for good and bad metrics"

$ws.Range("A5:F5").HorizontalAlignment = -4108
$ws.Range("A5:G5").VerticalAlignment = -4108
$ws.Range("A5:G5").WrapText = $true
$ws.Range("F5").Interior.Color = 5296274

# -----------------------------------------------------------------
# 6. Data row 6 (was row 5, nearly empty): day 3
# -----------------------------------------------------------------
$ws.Range("A6").Value = 3
$ws.Range("B6").Value = 43637
$ws.Range("B6").NumberFormat = "m/d/yyyy"

$ws.Range("A6:F6").HorizontalAlignment = -4108
$ws.Range("A6:G6").VerticalAlignment = -4108
$ws.Range("A6:G6").WrapText = $true
$ws.Range("C6").Font.Name = "Consolas"
$ws.Range("C6").Font.Size = 9
$ws.Range("C6").Font.Color = 5800553

# -----------------------------------------------------------------
# 8. Dates columns B4:B6 keep centered/wrapped with date format
# -----------------------------------------------------------------
$ws.Range("B4:B6").NumberFormat = "m/d/yyyy"
$ws.Range("B4:B6").HorizontalAlignment = -4108
$ws.Range("B4:B6").VerticalAlignment = -4108
$ws.Range("B4:B6").WrapText = $true

# -----------------------------------------------------------------
# 9. Column width adjustments
# -----------------------------------------------------------------
$ws.Columns.Item(6).ColumnWidth = 23.140625
$ws.Columns.Item(7).ColumnWidth = 15.85546875

# -----------------------------------------------------------------
# 10. Selection
# -----------------------------------------------------------------
$ws.Range("G6").Select()
